# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value (column G), rows 2..34
$newValues = @{
    2  = 5
    3  = 8
    4  = 5
    5  = 6
    6  = 6
    7  = 5
    8  = 2
    9  = 1
    10 = 4
    11 = 5
    12 = 11
    13 = 2
    14 = 5
    15 = 2
    16 = 4
    17 = 3
    18 = 8
    19 = 5
    20 = 4
    21 = 11
    22 = 5
    23 = 9
    24 = 4
    25 = 6
    26 = 3
    27 = 7
    28 = 5
    29 = 4
    30 = 12
    31 = 3
    32 = 6
    33 = 3
    34 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
